$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new September notification arrived ("bal axisbank" at 2024-09-23 08:35:16).
# It is the most recent entry, so it is inserted as a brand new row at the top
# of the September log (row 47), pushing every row from 47 down (through the
# "Others" data block and the trailing "Broadband" row) down by one row.
$ws.Rows("47:47").Insert()

$ws.Range("R47").Value = "bal axisbank"
$ws.Range("S47").Value = "2024-09-23 08:35:16"
